$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.792.22'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '3.414.45'
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.51'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -8.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.415.40'
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.482'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.50%  '
$ws.Range('E10').Value = '  -9.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.01'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -10.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.373'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -9.85%  '
$ws.Range('D13').Value = '3.992.39'
$ws.Range('E13').Value = '  -4.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000177'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -9.82%  '
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.09'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -10.45%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '64.675.05'
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.377.16'
$ws.Range('E18').Value = '  -5.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.38'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -15.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.68'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -8.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -8.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '380.23'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -10.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.539'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -9.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.70'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -7.78%  '
$ws.Range('D26').Value = '3.550.45'
$ws.Range('E26').Value = '  -4.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000104'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -11.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -11.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.18'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -12.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.94'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -12.15%  '
$ws.Range('D32').Value = '3.433.25'
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.143'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '22.83'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '170.23'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -13.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.60'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -14.41%  '
$ws.Range('E39').Value = '  -12.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.59'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -13.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0757'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -8.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.800'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.59%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.02'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.25'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -15.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.60'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -11.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.32'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -9.48%  '
$ws.Range('D50').Value = '2.197.40'
$ws.Range('E50').Value = '  -5.99%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.95'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -19.90%  '
